# Applies the "software engineering files" text cleanup described in the
# commit: a handful of adjacent same-formatted runs that PowerPoint had
# split apart (e.g. while the author was typing / autocorrect fired) are
# re-merged into single runs carrying the same text.
#
# Strategy: for each target paragraph we locate the desired final text as a
# contiguous substring of the shape's full TextRange.Text (dynamically, via
# IndexOf, so we don't depend on fragile hard-coded character offsets) and
# then re-assign that exact substring back onto the matching Characters()
# sub-range. PowerPoint's text engine collapses the runs spanned by that
# sub-range into a single run (taking on the formatting of the first run in
# the span), which is exactly the run-merge the diff shows.

function Merge-Text {
    param(
        $TextRange,
        [string]$Target
    )

    $full = $TextRange.Text
    $idx = $full.IndexOf($Target)
    if ($idx -lt 0) {
        # Already applied (or nothing to do) -- leave untouched.
        return
    }
    $start = $idx + 1
    $len = $Target.Length
    $sub = $TextRange.Characters($start, $len)
    $sub.Text = $Target
}

$p = $ppt.ActivePresentation

# --- Slide 2 : "Introduction" -------------------------------------------------
$slide2 = $p.Slides.Item(2)
$body2 = $slide2.Shapes.Item(2).TextFrame.TextRange

Merge-Text $body2 "3X3 "
Merge-Text $body2 "For this project the game will be change from its most common size 3x3 to a 6x6 with 4 squares require to win rather than 3."

# --- Slide 3 : "Purpose of the Game" ------------------------------------------
$slide3 = $p.Slides.Item(3)
$title3 = $slide3.Shapes.Item(1).TextFrame.TextRange

Merge-Text $title3 "urpose of the Game"

# --- Slide 5 : "User Requirement" ---------------------------------------------
$slide5 = $p.Slides.Item(5)
$body5 = $slide5.Shapes.Item(2).TextFrame.TextRange

Merge-Text $body5 "user requirement:"
